$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "verifyContactHeaders"

$ws.Range("A1").Value = "ContactHeader"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A3").Value = "Name"
$ws.Range("A4").Value = "Address"
$ws.Range("A5").Value = "Category"
$ws.Range("A6").Value = "Status"
$ws.Range("A7").Value = "Phone"
$ws.Range("A8").Value = "Email"
$ws.Range("A9").Value = "Options"

$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("H35").Select() | Out-Null
